# Comfenalco Cartagena - Estado de Cuenta
# - Adds a new "2509" period record (row) for worker ANGIE PAOLA MURILLO RODRIGUEZ (CC 1143410271)
# - Updates the accumulated "VALOR MORA" total and the "Cant. Periodos" counter accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Insert a new row right after the last worker record (row 20), pushing the
#     signature block (rows 25-26) down by one (to rows 26-27) ---
$ws.Rows.Item(21).Insert(-4121)

# Copy the formatting of the last existing worker row (20) onto the new row (21)
# so the new record keeps the same borders/fill/font/number-format as its neighbours.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Populate the new worker/period record ---
$ws.Range("B21").Value2 = "CC"
$ws.Range("C21").Value2 = "1143410271"
$ws.Range("D21").Value2 = "ANGIE PAOLA MURILLO RODRIGUEZ"
$ws.Range("E21").Value2 = "2509"
$ws.Range("F21").Value2 = 56940
$ws.Range("G21").Value2 = 1423500
$ws.Range("H21").Value2 = ""
$ws.Range("I21").Value2 = ""
$ws.Range("J21").Value2 = ""

# --- Update the summary figures to reflect the new period ---
# VALOR MORA (E11): old total 279006 + new period's 56940 = 335946
$ws.Range("E11").Value2 = 335946

# Cant. Periodos (F13): one more period now being reported (5 -> 6)
$ws.Range("F13").Value2 = 6

Write-Output "Edit applied"
